$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metric")
$ws.Activate()

# Update the "description" column (D) for the six Care Force metric rows to
# append the new "pp = percentage point" legend blurb.

$ws.Range("D6").Value = "<b>The Formal Care Force Participation</b> metric tracks the share of the total population (aged 18 and over) employed in formal care occupations among all those working or looking for work. <b>pp = percentage point</b> "

$ws.Range("D7").Value = "<b>The Informal Care Force Participation</b> metric tracks the share of the total population (aged 18 and over) engaged in at least 3 hours of informal, unpaid care work in a day. <b>pp = percentage point</b>"

$ws.Range("D8").Value = "<b>Formal Care Hours Worked</b> refers to the total hours worked in a day in paid care jobs by the entire U.S. population (aged 18 and over), also shown as the share of all formal hours worked in a day. <b>pp = percentage point</b> "

$ws.Range("D9").Value = "<b>Informal Care Hours Worked</b> refers to the total hours worked in unpaid care activities in a day by the entire U.S. population (aged 18 and over), also shown as the share of all informal hours worked in a day. <b>pp = percentage point</b> "

$ws.Range("D10").Value = "<b>Formal Economic Value of Care</b> estimates the total annual economic contribution of formal care work by aggregating the salaries of all formal care workers. <b>pp = percentage point</b> "

$ws.Range("D11").Value = "<b>Informal Economic Value of Care</b> estimates the total annual economic contribution of unpaid care work by multiplying the total hours spent in informal care activities by the federal minimum wage. This is a lower-bound estimate of the value of informal care. <b>pp = percentage point</b> "

# The longer descriptions now wrap onto more lines, so the rows housing them
# grow taller.
$ws.Rows.Item(7).RowHeight = 90
$ws.Rows.Item(8).RowHeight = 90
$ws.Rows.Item(9).RowHeight = 105
$ws.Rows.Item(11).RowHeight = 120

# Scroll the view down a bit and leave the selection on D6, matching where
# the editor ended up after making the change.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D6").Select()
